# Ke hoach lam viec (work plan update)
#
# The "Fix Bug" task (row 12) is now finished, so its STATUS cell (F12)
# moves from "In Progress" to "Completed" - the same status/format already
# used by every other finished task in the plan (e.g. F11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the status text itself.
$ws.Range("F12").Value = "Completed"

# Pick up the same look-and-feel as the other "Completed" rows by copying
# the format from F11 (an existing Completed cell) onto F12.
$ws.Range("F11").Copy() | Out-Null
$ws.Range("F12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# Reflect the user's final selection on the sheet (now spanning through the
# newly-updated F12 cell).
$ws.Range("F9:F12").Select() | Out-Null
